# Update the "2024" sheet: a new September entry ("balance your axis") was
# recorded for the "Others" group, pushing the existing September entries
# (broker / amazeloan x4) down by one row, and a new August "hdfc" entry was
# recorded as well, pushing the existing August entries (hdfc x4) down by one
# row and moving the "Broadband" group label from row 38 to the new row 39.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# --- September (Details/Date in columns R/S) for the "Others" group ---
# New entry inserted at the top; everything else shifts down by one row.
$ws.Range("R29").Value = "balance your axis"
$ws.Range("S29").Value = "2024-09-03 11:21:30"

$ws.Range("R30").Value = "broker"
$ws.Range("S30").Value = "2024-09-01 22:35:38"

$ws.Range("R31").Value = "amazeloan"
$ws.Range("S31").Value = "2024-09-01 10:12:03"

$ws.Range("R32").Value = "amazeloan"
$ws.Range("S32").Value = "2024-09-01 09:42:38"

$ws.Range("R33").Value = "amazeloan"
$ws.Range("S33").Value = "2024-09-01 09:29:24"

$ws.Range("R34").Value = "amazeloan"
$ws.Range("S34").Value = "2024-09-01 09:27:06"

# Row 34 previously held an August (P/Q) entry; that entry moves down too,
# so clear the old August values that used to live on row 34.
$ws.Range("P34").Value = ""
$ws.Range("Q34").Value = ""

# --- August (Details/Date in columns P/Q) for the "Others" group ---
# New entry inserted at the top of this sub-list (row 35); existing entries
# shift down by one row each.
$ws.Range("P35").Value = "hdfc"
$ws.Range("Q35").Value = "2024-08-30 12:15:48"

$ws.Range("P36").Value = "hdfc"
$ws.Range("Q36").Value = "2024-08-21 20:17:10"

$ws.Range("P37").Value = "hdfc"
$ws.Range("Q37").Value = "2024-08-21 20:16:45"

$ws.Range("P38").Value = "hdfc"
$ws.Range("Q38").Value = "2024-08-21 20:15:50"

# The "Broadband" group label used to be on row 38; it now moves to the new
# row 39 since the "Others" group grew by one row.
$ws.Range("A38").Value = ""
$ws.Range("A39").Value = "Broadband"
